$wb = $excel.ActiveWorkbook

# Existing "perf stat" worksheet is currently the active tab; the new
# "papi stat" worksheet (PAPI L3 cache-miss counters) is inserted right
# after it, mirroring the two new "loop body papi" / "barrier papi"
# chart sheets that were added alongside it upstream.
$perfStat = $wb.Worksheets.Item("perf stat")

$papiStat = $wb.Worksheets.Add($null, $perfStat)
$papiStat.Name = "papi stat"

# Headers
$papiStat.Range("A1").Value = "Threads"
$papiStat.Range("B1").Value = "L3 Cache Misses: for (loop body)"
$papiStat.Range("C1").Value = "L3 Cache Misses: for (barrier enter/exit)"

# Data rows: thread count, loop-body L3 misses, barrier L3 misses
$papiStat.Range("A2").Value = 1
$papiStat.Range("B2").Value = 243000000
$papiStat.Range("C2").Value = 76050

$papiStat.Range("A3").Value = 2
$papiStat.Range("B3").Value = 147300000
$papiStat.Range("C3").Value = 59330

$papiStat.Range("A4").Value = 4
$papiStat.Range("B4").Value = 68410000
$papiStat.Range("C4").Value = 43600

$papiStat.Range("A5").Value = 8
$papiStat.Range("B5").Value = 35990000
$papiStat.Range("C5").Value = 32180

# Column widths to fit the long header text (matches the "perf stat"-style sheets)
$papiStat.Columns.Item(1).ColumnWidth = 10.83
$papiStat.Columns.Item(2).ColumnWidth = 36.33
$papiStat.Columns.Item(3).ColumnWidth = 46.33

# Numeric columns use a "0.00" number format, headers are bold-ish like the
# other summary sheets
$papiStat.Range("B2:C5").NumberFormat = "0.00"
$papiStat.Range("A1:C1").Font.Bold = $true
$papiStat.Range("A1:A5").HorizontalAlignment = -4108

# "papi stat" becomes the newly selected/active sheet (mirrors tabSelected
# moving from "perf stat" to "papi stat"), with B2 highlighted.
$papiStat.Range("B2").Select()
$papiStat.Activate()
